$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "ECs" sending-cluster row (original row 2). This shifts the
# "FAPs" row up to row 2 and the "MuSCs" row up to row 3, matching the
# new layout (only 2 data rows remain).
$ws.Rows("2:2").Delete()

# Update the recomputed TPM-derived values for the (now) row 2 - "FAPs".
$ws.Range("G2").Value = 0.115543
$ws.Range("H2").Value = 0.346629
$ws.Range("I2").Value = 0.7111155332715143
$ws.Range("J2").Value = 0.7111155332715143
$ws.Range("Q2").Value = 0.03027457686
$ws.Range("R2").Value = 0.27247119174
$ws.Range("S2").Value = 0.7111155332715143
$ws.Range("T2").Value = 0.7111155332715143

# Update the recomputed TPM-derived values for the (now) row 3 - "MuSCs".
$ws.Range("G3").Value = 0.04693833333333333
$ws.Range("H3").Value = 0.140815
$ws.Range("I3").Value = 0.2888844667284857
$ws.Range("J3").Value = 0.2888844667284857
$ws.Range("Q3").Value = 0.0122987821
$ws.Range("R3").Value = 0.1106890389
$ws.Range("S3").Value = 0.2888844667284857
$ws.Range("T3").Value = 0.2888844667284857
